$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (dbExcel) and E (WebExcel) file names for rows 2-4
$dFile = @'
TC13_CDS_Filter_InstrumentModel-NotSpecifiedindata_Neo4jData.xlsx
'@
$eFile = @'
TC13_CDS_Filter_InstrumentModel-NotSpecifiedindata_WebData.xlsx
'@
$ws.Range("D2").Value = $dFile
$ws.Range("D3").Value = $dFile
$ws.Range("D4").Value = $dFile
$ws.Range("E2").Value = $eFile
$ws.Range("E3").Value = $eFile
$ws.Range("E4").Value = $eFile

# Column B (query) per row, instrument model filter updated
$participantsQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Not specified in data']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p, s, collect(distinct samp.sample_id) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY `Participant ID`LIMIT 100
'@
$samplesQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Not specified in data']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
ORDER By samp.sample_id LIMIT 100
'@
$filesQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Not specified in data']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,f,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN 
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER By f.file_name LIMIT 100
'@
$ws.Range("B2").Value = $participantsQuery
$ws.Range("B3").Value = $samplesQuery
$ws.Range("B4").Value = $filesQuery

# Column C (StatQuery) - same for all 3 rows, instrument model filter updated
$statQuery = @'
MATCH (f:file)
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Not specified in data']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,f, s, collect(distinct samp.sample_id) as samp
RETURN
count(distinct s) AS Studies,
count(distinct p) AS Participants,
count(distinct samp) AS Samples,
count(distinct f) AS Files
'@
$ws.Range("C2").Value = $statQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("C4").Value = $statQuery

# Re-fit columns D and E to their new (longer) best-fit widths,
# matching the 'bestFit' recompute Excel performs after the content grows
$ws.Columns.Item(4).ColumnWidth = 93.66666666666667
$ws.Columns.Item(5).ColumnWidth = 92.16666666666667

# Update selection to D2 as last active cell
$ws.Range("D2").Select()
